# Rework the small item table:
#  - rename "number" -> "quantity", "rating" -> "price", and add a new "total" column
#  - restyle the header row and the existing data columns (this is what produces the
#    extra cell style record that appears in the saved workbook)
#  - leave the new "price"/"total" data cells for rows 2-6 blank (only the headers are added)
#  - finish with the active cell on D2, under the new "total" header

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update/extend the header row: B1 "number" -> "quantity", C1 "rating" -> "price",
# and the brand new D1 "total" header.
$ws.Range("B1").Value = "quantity"
$ws.Range("C1").Value = "price"
$ws.Range("D1").Value = "total"

# Re-apply the base cell style across the full header row (A1:D1) and across the
# existing item/quantity data columns (A2:B6).
$ws.Range("A1:D1").Style = "Normal"
$ws.Range("A2:B6").Style = "Normal"

# Match the workbook's final selection state.
$ws.Range("D2").Select() | Out-Null
